# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" (right after "总计") holding the newly
# reported fund-holding data for that quarter, and updates the "总计"
# (summary) sheet with the corresponding new row. The pre-existing
# "2022-Q2", "2022-Q1" and "2021-Q3" sheets are left as-is (they simply
# shift one tab to the right to make room for the new quarter).

$wb = $excel.ActiveWorkbook

# --- locate the sheet that used to sit right after "总计" (currently "2022-Q2") ---
$refSheet = $wb.Worksheets.Item(2)

# --- insert the brand-new quarter sheet just before it ---
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Columns B and D:G hold text-formatted values (fund code keeps its
# leading zero, and the numeric-looking figures are stored as text),
# matching the layout used by the other quarterly sheets.
$q3.Range("B1:B3").NumberFormat = "@"
$q3.Range("D1:G3").NumberFormat = "@"

# Header row (same layout/wording as the other quarterly sheets)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row 2 - Class A units
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "011243"
$q3.Range("C2").Value = "万家惠裕回报6个月持有期混合A"
$q3.Range("D2").Value = "1.54"
$q3.Range("E2").Value = "27.67"
$q3.Range("F2").Value = "0.76"
$q3.Range("G2").Value = "0.0117"
$q3.Range("H2").Value = 7

# Row 3 - Class C units
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "011244"
$q3.Range("C3").Value = "万家惠裕回报6个月持有期混合C"
$q3.Range("D3").Value = "0.12"
$q3.Range("E3").Value = "27.67"
$q3.Range("F3").Value = "0.76"
$q3.Range("G3").Value = "0.0009"
$q3.Range("H3").Value = 7

# Match the bold/centered/boxed look used for the header row and the
# running-index column on the sibling quarterly sheets.
$refSheet.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$refSheet.Range("A2").Copy()
$q3.Range("A2:A3").PasteSpecial(-4122)

# --- update the "总计" summary sheet: insert a new top row for 2022-Q3 ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Renumber the running index in column A for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Fill in the new 2022-Q3 row
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

# The new row's data cells (B2:D2) should stay plain, unstyled cells like
# the other data rows; only the index column carries the boxed style.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Keep "2021-Q3" as the selected/active tab, same as before the edit.
$wb.Worksheets.Item("2021-Q3").Activate()
